# Moved Video to useage bullets, as they were breaking in the screenshots section.
#
# The "Video Demo of Application" block (two rows: a label row and a URL
# row, tagged VIDEO) is removed from the Screenshots area and folded into
# a single bullet ("<label> - <url>") appended to the "Use" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the two existing Video cells (label + url) before we touch anything.
$videoLabel = $ws.Range("B23").Value2
$videoUrl   = $ws.Range("B24").Value2
$videoText  = $videoLabel + " - " + $videoUrl

# Remove the whole Video block (rows 23 and 24) - this also removes the
# now-orphaned "Video Demo of Application" / url / "VIDEO" shared strings.
$ws.Rows("23:24").Delete()

# Insert a new row right after the last "Use" bullet (old row 17) to hold
# the merged Video bullet.
$ws.Rows("18:18").Insert()

$ws.Range("B18").Value = $videoText
$ws.Range("C18").Value = "Use"

# Match the author's final selection.
$ws.Range("H22").Select() | Out-Null
